$wb = $excel.ActiveWorkbook

# Work on Sheet2: add new trial data
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("D2").Value = 336
$ws2.Range("B4").Formula = "=0.1411*`$D`$2-1.7937"

# Update selection on Sheet2 to D3
$ws2.Range("D3").Select()

# Make Sheet2 the active sheet/tab (activeTab="1" -> second sheet, index 2)
$ws2.Activate()

$wb.Save()
